$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new data row ("2026/01/30", "金", 23, 201) was inserted right before the
# existing row 752 ("2026/12/29", ...), pushing all following rows down by
# one (old row 752 -> new row 753, ..., old row 793 -> new row 794).
$ws.Rows.Item(752).Insert()

# Column A holds dates as literal text (e.g. "2026/01/30"), not real Excel
# dates. A leading apostrophe forces the COM layer to store it as text
# instead of auto-converting it to a date serial number; resetting the
# style back to "Normal" afterwards clears the quote-prefix/number-format
# that the text coercion leaves behind on the cell.
$ws.Range("A752").Value = "'2026/01/30"
$ws.Range("B752").Value = "金"
$ws.Range("C752").Value = 23
$ws.Range("D752").Value = 201
$ws.Range("A752:D752").Style = "Normal"
